$d = $word.ActiveDocument
$inf = [char]0x221E

# --- Change 1: paragraph 11) "... x = 7" -> append "; (x = 4 is extraneous)" as a new run
$p11 = $d.Paragraphs.Item(7)
$r11 = $p11.Range
$r11.SetRange($r11.End - 1, $r11.End - 1)
$r11.InsertAfter(";  (x = 4 is extraneous)")

# --- Change 2: paragraph 20) remove the gramStart/gramEnd proofErr markers and
#     merge " ," + " 12]" into a single run " , 12]", while preserving the
#     Cambria Math formatted infinity run and keeping "20)" / " (-" as separate runs.

# Step 2a: merge " ," and " 12]" into one run " , 12]" (this also removes the
# now-enclosed gramEnd proofErr marker).
$null = $d.Content.Find.Execute(" , 12]", $true, $false, $false, $false, $false, $true, 1, $false, " , 12]", 2)

# Step 2b: remove the gramStart proofErr marker. Replacing "-<inf>" merges the
# surrounding plain runs (including "20)") into one, and drops the Cambria Math
# formatting on the infinity character, so we restore both afterwards.
$null = $d.Content.Find.Execute("-$inf", $true, $false, $false, $false, $false, $true, 1, $false, "-$inf", 2)

# Step 2c: restore the Cambria Math formatting on the infinity character, which
# also splits it back out into its own run.
$rngInf = $d.Content
$null = $rngInf.Find.Execute($inf, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngInf.Font.Name = "Cambria Math"

# Step 2d: split "20)" back out into its own run (it got merged with " (-" in
# step 2b) by clearing and re-inserting its text as plain literal text.
$rng20 = $d.Content
$null = $rng20.Find.Execute("20)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng20.Text = ""
$rng20.InsertBefore("20)")

# --- Change 3: paragraph 26) -> append " [-1, 5]" as a new run.
# A plain InsertAfter would merge into the existing run since neither run
# carries distinguishing metadata, so we use track-changes insert + accept of
# just that revision, which preserves the run boundary without side effects.
$p26 = $d.Paragraphs.Item($d.Paragraphs.Count)
$d.TrackRevisions = $true
$r26 = $p26.Range
$r26.SetRange($r26.End - 1, $r26.End - 1)
$r26.InsertAfter(" [-1, 5]")
$d.TrackRevisions = $false
$d.Revisions.Item(1).Accept()
